# Daily attendance processing - 2026-01-12 19:12:20
# Normalize the "Recorded By" (column G) entries by reversing the order
# of the comma-separated recorder list for specific known values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -eq "backup@backdoor.com, System, system") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
    elseif ($current -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($current -eq "backup@backdoor.com, System") {
        $cell.Value = "System, backup@backdoor.com"
    }
}
